$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure B/C columns for rows 3-6 are written as text (matching inlineStr/t="s" in source)
# even though their contents look like numbers. Temporarily force a text number
# format so Excel doesn't silently convert the value to a real number, then
# restore the original (Normal) style so no stray style index is left behind.
$ws.Range("B3:C6").NumberFormat = "@"

# Row 3 (flow_scour)
$ws.Range("B3").Value = "5"
$ws.Range("C3").Value = "0"
$ws.Range("D3").Value = 0.04
$ws.Range("E3").Value = 0

# Row 4 (flow_base_flow)
$ws.Range("B4").Value = "90"
$ws.Range("C4").Value = "0"
$ws.Range("D4").Value = 0.71
$ws.Range("E4").Value = 0

# Row 5 (food_web)
$ws.Range("B5").Value = "5"
$ws.Range("C5").Value = "30"
$ws.Range("D5").Value = 0.04
$ws.Range("E5").Value = 0.24

# Row 6 (temperature_rearing)
$ws.Range("B6").Value = "54"
$ws.Range("C6").Value = "0"
$ws.Range("D6").Value = 0.43
$ws.Range("E6").Value = 0

# Restore the default (unstyled) cell style now that the values are locked in as text.
$ws.Range("B3:C6").Style = "Normal"
